$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the Date of Birth value in D2 (master/secondary extracts no longer
# always carry a DOB) and drop its "Date" formatting back to the default.
$ws.Range("D2").ClearContents()
$ws.Range("D2").Style = "Normal"

# The "Date" cell style is no longer used anywhere in the workbook -
# remove it from the style gallery.
$wb.Styles.Item("Date").Delete()

# Column D no longer needs to be as wide now that it holds no dates.
$ws.Columns.Item(4).ColumnWidth = 15.17
